# Add season record columns (Wins, Losses, Ties) to the player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells AD1:AF1 - copy formatting from the existing header style (AC1)
# then overwrite with the new header text.
$headerStyleSource = $ws.Range("AC1")

$wins = $ws.Range("AD1")
$headerStyleSource.Copy($wins)
$wins.Value = "Wins"

$losses = $ws.Range("AE1")
$headerStyleSource.Copy($losses)
$losses.Value = "Losses"

$ties = $ws.Range("AF1")
$headerStyleSource.Copy($ties)
$ties.Value = "Ties"

# Data rows 2:41 - season record values (same for every row in this sheet)
$ws.Range("AD2:AD41").Value = 80
$ws.Range("AE2:AE41").Value = 81
$ws.Range("AF2:AF41").Value = 0
